$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update the data values in row 4 (E4, G4, H4)
$ws.Range("E4").Value = 6
$ws.Range("G4").Value = 3
$ws.Range("H4").Value = 13

# Move the active selection to A4 (was I7)
$ws.Range("A4").Select()
